# Update countries & provincias Spain
# Refresh COVID data snapshot: update "last updated" timestamp and refresh
# several country rows with new case counts. Because the sheet is kept
# sorted by "Casos totales" (column B) descending, a few rows swap order
# with their neighbour as a result of the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$row, [string]$pais, [int]$casosTotales, [int]$nuevosCasos, [int]$casosActivos, [int]$recuperados, [int]$casosCriticos, [int]$muertesHoy, [int]$muertes) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- Header: timestamp of last update -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 11:52"

# --- Rows 37-39: Indonesia overtakes Dinamarca and Bielorrusia ------------
Set-Row 37 "Indonesia"   8211 436 1002 6520  0 42 689
Set-Row 38 "Dinamarca"   8210 137 5384 2432 74  0 394
Set-Row 39 "Bielorrusia" 8022   0  938 7024 92  0  60

# --- Row 52: Finlandia updated counts --------------------------------------
Set-Row 52 "Finlandia" 4395 111 2000 2223 60 0 172

# --- Row 57: Argentina updated counts --------------------------------------
Set-Row 57 "Argentina" 3435 0 919 2350 123 1 166

# --- Rows 71-72: Estonia overtakes Armenia ---------------------------------
Set-Row 71 "Estonia" 1605 13  206 1353  6 1 46
Set-Row 72 "Armenia" 1596 73  728  841 10 3 27

# --- Rows 98-99: Albania overtakes Niger ------------------------------------
Set-Row 98 "Albania" 678 15 394 257 4 0 27
Set-Row 99 "Niger"   671  0 256 391 0 0 24

# --- Row 140: Etiopia updated counts ----------------------------------------
Set-Row 140 "Etiopia" 117 1 25 89 0 0 3
